$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "BAU Emissions" sheet: relabel the " : NoSettings" suffix to " : test"
#    across all of the row labels in column A, update the forecast values in
#    row 94 (columns M:AE), and move the on-screen selection.
# ---------------------------------------------------------------------------
$wsBau = $wb.Worksheets.Item("BAU Emissions")
$wsBau.Activate()

$wsBau.Range("A4:A280").Replace(" : NoSettings", " : test")

$wsBau.Range("M94").Value = 1001080
$wsBau.Range("N94").Value = 2002150
$wsBau.Range("O94").Value = 3003230
$wsBau.Range("P94").Value = 4004300
$wsBau.Range("Q94:AE94").Value = 5005380

$wsBau.Range("A30:AE280").Select()

# ---------------------------------------------------------------------------
# 2. "About" sheet: bump the last-updated date and make this the active tab.
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("C1").Value = 45387
$wsAbout.Range("E29").Select()
